$d = $word.ActiveDocument

# --- 1. Remove the stray "_GoBack" bookmark from the empty paragraph that
#        immediately precedes the "Ultima Online (house decay)" heading ---
$headingRange = $d.Content
$found = $headingRange.Find.Execute("Ultima Online (house decay)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $headingPara = $headingRange.Paragraphs(1)
    $bmPara = $d.Paragraphs($headingPara.Index - 1)
    $bmXml = "<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr></w:p>"
    $bmPara.Range.InsertXML($bmXml)
}

# --- 2. Append the new "Guns of Icarus" research section at the end of
#        the document ---
$endRange = $d.Content
$endRange.Collapse(0)
$newBlock = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:rPr>
          <w:color w:val="7030A0"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="7030A0"/>
        </w:rPr>
        <w:t>Further research:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:sz w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>Guns of Icarus</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Gameplay is very similar to our design (close to a first-person experience of groups current design)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Cooperative gameplay is a huge element of the game.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Players must coordinate to counter attackers, while also prioritising how to manage incoming damage and the restrictions this causes.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Player </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>positioning</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> and skill level are a significant part of the game</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Positive feedback loops – if a team’s ship becomes damaged, becomes increasingly harder to survive</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Game creates very intense state of panic throughout the team – all team members must be communicating while actively completing tasks and predicting future tasks</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Game features </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>a player class system where certain players</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> are more adept at certain tasks (complete faster, use less resources) and only able to complete specific tasks</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Some items exist as tools within the level and must be picked up and taken to complete the task at a component</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> – add functionality that is not present in player classes</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="6"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Repairing and mitigating damage is done is stages. Rather </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">than a set </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>time period</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> from when players begin a task, each player interaction restores a portion of the activities capacity. This allows for higher pressure situations where performing perfectly will only maintain – at best – the current situation.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
     </w:p>
'@
$endRange.InsertXML($newBlock)
